# Updating the hazard IDs
# The table used to start at row 2 (header) with data rows 3-6, followed by
# several trailing empty formatted rows (7-14). The new layout moves the
# table up so the header is row 1 and the data is rows 2-5, the trailing
# empty rows are removed, and the "ID" column now holds text IDs (H-1..H-4)
# instead of plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-empty leading row so the header moves from row 2 to row 1
# (and everything below shifts up with it).
$ws.Rows("1:1").Delete()

# Remove the trailing empty rows (now at rows 6:13 after the shift above).
$ws.Rows("6:13").Delete()

# Replace the numeric hazard IDs with textual hazard IDs.
$ws.Range("A2").Value = "H-1"
$ws.Range("A3").Value = "H-2"
$ws.Range("A4").Value = "H-3"
$ws.Range("A5").Value = "H-4"

# The last data row's content got shorter/reflowed -> smaller row height.
$ws.Rows("5:5").RowHeight = 75

# Probable-causes cells for the 3rd/4th data rows switch from
# top-aligned to vertically-centered wrapped text (matching the other rows).
$ws.Range("C4").VerticalAlignment = -4108
$ws.Range("D4").VerticalAlignment = -4108
$ws.Range("D5").VerticalAlignment = -4108

# Column D widened to fit the new content, no longer relying on auto bestFit.
$ws.Columns("D").ColumnWidth = 42.67

# Selection moves to A2 (first data row) instead of B3.
$ws.Range("A2").Select()
